# "add data and update plotting"
#
# Renames two headers/labels (space -> underscore) and regenerates the
# L column (mut freq) formulas as a filled-down shared-formula block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename labels: "TAD1 C184A" -> "TAD1_C184A", "mut freq" -> "mut_freq" ---
$ws.Range("C5").Value = "TAD1_C184A"
$ws.Range("C6").Value = "TAD1_C184A"
$ws.Range("C7").Value = "TAD1_C184A"
$ws.Range("C11").Value = "TAD1_C184A"
$ws.Range("C12").Value = "TAD1_C184A"
$ws.Range("C13").Value = "TAD1_C184A"

# Match style of the already-renamed rows (11-13) by copying style to 5-7
$ws.Range("C11").Copy()
$ws.Range("C5:C7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("L1").Value = "mut_freq"

# --- Re-enter the mut-freq formula in L2 and fill it down through L7 so
#     Excel regenerates it as one shared formula group ---
$ws.Range("L2").Formula = "=K8/K2"
$ws.Range("L2:L7").Formula = "=K8/K2"

# --- Update the view: scroll so column F is the left-most visible column
#     and select L1 ---
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("L1").Select()
